$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.911.76'
$ws.Range('E2').Value = '  -1.74%  '
$ws.Range('D3').Value = '1.972.82'
$ws.Range('E3').Value = '  -3.54%  '
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '238.42'
$ws.Range('E5').Value = '  -7.29%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.599'
$ws.Range('E6').Value = '  -4.19%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '53.73'
$ws.Range('E8').Value = '  -6.64%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '59.78'
$ws.Range('E9').Value = '  +4.64%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.368'
$ws.Range('E10').Value = '  -5.39%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0748'
$ws.Range('E11').Value = '  -6.59%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0982'
$ws.Range('E12').Value = '  -4.82%  '
$ws.Range('D13').Value = '2.265.39'
$ws.Range('E13').Value = '  -3.47%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '13.82'
$ws.Range('E14').Value = '  -6.93%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '20.62'
$ws.Range('E15').Value = '  -3.86%  '
$ws.Range('E16').Value = '  -9.23%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.00'
$ws.Range('E17').Value = '  -7.33%  '
$ws.Range('D18').Value = '1.974.75'
$ws.Range('E18').Value = '  -3.46%  '
$ws.Range('D19').Value = '36.848.25'
$ws.Range('E19').Value = '  -1.76%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '68.00'
$ws.Range('D21').Value = '0.0₃0804'
$ws.Range('E21').Value = '  -6.17%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '226.57'
$ws.Range('E22').Value = '  -1.14%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.93'
$ws.Range('E23').Value = '  -6.25%  '
$ws.Range('E24').Value = '  +0.02%  '
$ws.Range('E25').Value = '  -0.39%  '
$ws.Range('E26').Value = '  -12.26%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '160.91'
$ws.Range('E27').Value = '  -1.65%  '
$ws.Range('E28').Value = '  -6.76%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.93'
$ws.Range('E29').Value = '  -5.18%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.123'
$ws.Range('E30').Value = '  -11.29%  '
$ws.Range('E31').Value = '  -5.78%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.117'
$ws.Range('E32').Value = '  -3.48%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.38'
$ws.Range('E33').Value = '  -8.24%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0610'
$ws.Range('E34').Value = '  -8.73%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.21'
$ws.Range('E35').Value = '  -7.21%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.30'
$ws.Range('E36').Value = '  -8.72%  '
$ws.Range('E37').Value = '  -0.10%  '
$ws.Range('E38').Value = '  -1.46%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.27'
$ws.Range('E39').Value = '  -5.74%  '
$ws.Range('E40').Value = '  -4.39%  '
$ws.Range('D42').Value = '1.412.35'
$ws.Range('E42').Value = '  +0.05%  '
$ws.Range('E43').Value = '  -6.39%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0886'
$ws.Range('E44').Value = '  -8.48%  '
$ws.Range('E45').Value = '  -7.47%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '86.99'
$ws.Range('E46').Value = '  -4.98%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '15.10'
$ws.Range('E47').Value = '  -7.38%  '
$ws.Range('B48').Value = 'ARBITRUM'
$ws.Range('C48').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.987'
$ws.Range('E48').Value = '  -6.46%  '
$ws.Range('B49').Value = 'MXToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.85'
$ws.Range('E49').Value = '  -0.94%  '
$ws.Range('B50').Value = 'FTXToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.72'
$ws.Range('E50').Value = '  +13.81%  '
$ws.Range('E51').Value = '  -11.65%  '
